$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Duration_Seconds" values for the 4 trial rows (3, 5, 7, 9) from 1 -> 10
$ws.Range("C3").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("C9").Value = 10

# Row 3 (Trial 1): condition stays "3D", filenames change to Face1_L.jpg / Face1_R.jpg
$ws.Range("B3").Value = "3D"
$ws.Range("D3").Value = "Face1_L.jpg"
$ws.Range("E3").Value = "Face1_R.jpg"

# Row 5 (Trial 2): condition stays "MonocL", left filename changes to Face1_L.jpg
$ws.Range("B5").Value = "MonocL"
$ws.Range("D5").Value = "Face1_L.jpg"

# Row 7 (Trial 3): condition changes from "MonocR" to "Pseudo", filenames set
$ws.Range("B7").Value = "Pseudo"
$ws.Range("D7").Value = "Face1_R.jpg"
$ws.Range("E7").Value = "Face1_L.jpg"

# Row 9 (Trial 4): condition stays "2D", filenames change to Face1_R.jpg / Face1_R.jpg
$ws.Range("B9").Value = "2D"
$ws.Range("D9").Value = "Face1_R.jpg"
$ws.Range("E9").Value = "Face1_R.jpg"

$wb.Save()
